$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fix rows 503-507 (Name-column value belongs in column E, not F) ---
foreach ($r in 503..507) {
    $val = $ws.Cells.Item($r, 6).Value()
    $ws.Cells.Item($r, 4).Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4122)
    $ws.Cells.Item($r, 5).Value = $val
    $ws.Cells.Item($r, 6).Clear()
}
$excel.CutCopyMode = 0

# --- Step 2: append new rows 508-527, extracted structural data from 10.1016/j.intermet.2021.107167 ---
# Row 508: CoCrFeNi / hardness
$ws.Cells.Item(488, 2).Copy()
$ws.Cells.Item(508, 2).PasteSpecial(-4122)
$ws.Cells.Item(508, 2).Value = 'CoCrFeNi'
$ws.Cells.Item(508, 3).Value = 'FCC'
$ws.Cells.Item(508, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(508, 4).Copy()
$ws.Cells.Item(508, 5).PasteSpecial(-4122)
$ws.Cells.Item(508, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(508, 6).Value = 'hardness'
$ws.Cells.Item(508, 7).Value = 'EXP'
$ws.Cells.Item(508, 9).Value = 298
$ws.Cells.Item(488, 10).Copy()
$ws.Cells.Item(508, 10).PasteSpecial(-4122)
$ws.Cells.Item(488, 10).Copy()
$ws.Cells.Item(508, 11).PasteSpecial(-4122)
$ws.Cells.Item(508, 10).Formula = '=P508*9807000'
$ws.Cells.Item(508, 11).Formula = '=Q508*9807000'
$ws.Cells.Item(508, 12).Value = 'Pa'
$ws.Cells.Item(508, 13).Value = 'T2'
$ws.Cells.Item(508, 14).Value = '10.1016/j.intermet.2021.107167'
$ws.Cells.Item(508, 16).Value = 110
$ws.Cells.Item(508, 17).Value = 3

# Row 509: CoCrFeNiGe0.1 / hardness
$ws.Cells.Item(509, 2).Value = 'CoCrFeNiGe0.1'
$ws.Cells.Item(509, 3).Value = 'FCC'
$ws.Cells.Item(509, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(509, 4).Copy()
$ws.Cells.Item(509, 5).PasteSpecial(-4122)
$ws.Cells.Item(509, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(509, 6).Value = 'hardness'
$ws.Cells.Item(509, 7).Value = 'EXP'
$ws.Cells.Item(509, 9).Value = 298
$ws.Cells.Item(488, 10).Copy()
$ws.Cells.Item(509, 10).PasteSpecial(-4122)
$ws.Cells.Item(488, 10).Copy()
$ws.Cells.Item(509, 11).PasteSpecial(-4122)
$ws.Cells.Item(509, 10).Formula = '=P509*9807000'
$ws.Cells.Item(509, 11).Formula = '=Q509*9807000'
$ws.Cells.Item(509, 12).Value = 'Pa'
$ws.Cells.Item(509, 13).Value = 'T2'
$ws.Cells.Item(509, 14).Value = '10.1016/j.intermet.2021.107167'
$ws.Cells.Item(509, 16).Value = 113
$ws.Cells.Item(509, 17).Value = 3

# Row 510: CoCrFeNiGe0.2 / hardness
$ws.Cells.Item(510, 2).Value = 'CoCrFeNiGe0.2'
$ws.Cells.Item(510, 3).Value = 'FCC'
$ws.Cells.Item(510, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(510, 4).Copy()
$ws.Cells.Item(510, 5).PasteSpecial(-4122)
$ws.Cells.Item(510, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(510, 6).Value = 'hardness'
$ws.Cells.Item(510, 7).Value = 'EXP'
$ws.Cells.Item(510, 9).Value = 298
$ws.Cells.Item(488, 10).Copy()
$ws.Cells.Item(510, 10).PasteSpecial(-4122)
$ws.Cells.Item(488, 10).Copy()
$ws.Cells.Item(510, 11).PasteSpecial(-4122)
$ws.Cells.Item(510, 10).Formula = '=P510*9807000'
$ws.Cells.Item(510, 11).Formula = '=Q510*9807000'
$ws.Cells.Item(510, 12).Value = 'Pa'
$ws.Cells.Item(510, 13).Value = 'T2'
$ws.Cells.Item(510, 14).Value = '10.1016/j.intermet.2021.107167'
$ws.Cells.Item(510, 16).Value = 118
$ws.Cells.Item(510, 17).Value = 3

# Row 511: CoCrFeNiGe0.3 / hardness
$ws.Cells.Item(511, 2).Value = 'CoCrFeNiGe0.3'
$ws.Cells.Item(511, 3).Value = 'FCC'
$ws.Cells.Item(511, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(511, 4).Copy()
$ws.Cells.Item(511, 5).PasteSpecial(-4122)
$ws.Cells.Item(511, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(511, 6).Value = 'hardness'
$ws.Cells.Item(511, 7).Value = 'EXP'
$ws.Cells.Item(511, 9).Value = 298
$ws.Cells.Item(488, 10).Copy()
$ws.Cells.Item(511, 10).PasteSpecial(-4122)
$ws.Cells.Item(488, 10).Copy()
$ws.Cells.Item(511, 11).PasteSpecial(-4122)
$ws.Cells.Item(511, 10).Formula = '=P511*9807000'
$ws.Cells.Item(511, 11).Formula = '=Q511*9807000'
$ws.Cells.Item(511, 12).Value = 'Pa'
$ws.Cells.Item(511, 13).Value = 'T2'
$ws.Cells.Item(511, 14).Value = '10.1016/j.intermet.2021.107167'
$ws.Cells.Item(511, 16).Value = 119
$ws.Cells.Item(511, 17).Value = 2

# Row 512: CoCrFeNi / youngs modulus
$ws.Cells.Item(488, 2).Copy()
$ws.Cells.Item(512, 2).PasteSpecial(-4122)
$ws.Cells.Item(512, 2).Value = 'CoCrFeNi'
$ws.Cells.Item(512, 3).Value = 'FCC'
$ws.Cells.Item(512, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(512, 4).Copy()
$ws.Cells.Item(512, 5).PasteSpecial(-4122)
$ws.Cells.Item(512, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(512, 6).Value = 'youngs modulus'
$ws.Cells.Item(512, 7).Value = 'EXP'
$ws.Cells.Item(512, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(512, 9).Value = 298
$ws.Cells.Item(512, 10).Value = 214000000000
$ws.Cells.Item(512, 11).Value = 11000000000
$ws.Cells.Item(512, 12).Value = 'Pa'
$ws.Cells.Item(512, 13).Value = 'T2'
$ws.Cells.Item(512, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 513: CoCrFeNiGe0.1 / youngs modulus
$ws.Cells.Item(513, 2).Value = 'CoCrFeNiGe0.1'
$ws.Cells.Item(513, 3).Value = 'FCC'
$ws.Cells.Item(513, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(513, 4).Copy()
$ws.Cells.Item(513, 5).PasteSpecial(-4122)
$ws.Cells.Item(513, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(513, 6).Value = 'youngs modulus'
$ws.Cells.Item(513, 7).Value = 'EXP'
$ws.Cells.Item(513, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(513, 9).Value = 298
$ws.Cells.Item(513, 10).Value = 217000000000
$ws.Cells.Item(513, 11).Value = 8000000000
$ws.Cells.Item(513, 12).Value = 'Pa'
$ws.Cells.Item(513, 13).Value = 'T2'
$ws.Cells.Item(513, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 514: CoCrFeNiGe0.2 / youngs modulus
$ws.Cells.Item(514, 2).Value = 'CoCrFeNiGe0.2'
$ws.Cells.Item(514, 3).Value = 'FCC'
$ws.Cells.Item(514, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(514, 4).Copy()
$ws.Cells.Item(514, 5).PasteSpecial(-4122)
$ws.Cells.Item(514, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(514, 6).Value = 'youngs modulus'
$ws.Cells.Item(514, 7).Value = 'EXP'
$ws.Cells.Item(514, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(514, 9).Value = 298
$ws.Cells.Item(514, 10).Value = 218000000000
$ws.Cells.Item(514, 11).Value = 10000000000
$ws.Cells.Item(514, 12).Value = 'Pa'
$ws.Cells.Item(514, 13).Value = 'T2'
$ws.Cells.Item(514, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 515: CoCrFeNiGe0.3 / youngs modulus
$ws.Cells.Item(515, 2).Value = 'CoCrFeNiGe0.3'
$ws.Cells.Item(515, 3).Value = 'FCC'
$ws.Cells.Item(515, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(515, 4).Copy()
$ws.Cells.Item(515, 5).PasteSpecial(-4122)
$ws.Cells.Item(515, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(515, 6).Value = 'youngs modulus'
$ws.Cells.Item(515, 7).Value = 'EXP'
$ws.Cells.Item(515, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(515, 9).Value = 298
$ws.Cells.Item(515, 10).Value = 217000000000
$ws.Cells.Item(515, 11).Value = 13000000000
$ws.Cells.Item(515, 12).Value = 'Pa'
$ws.Cells.Item(515, 13).Value = 'T2'
$ws.Cells.Item(515, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 516: CoCrFeNi / tensile yield stress
$ws.Cells.Item(488, 2).Copy()
$ws.Cells.Item(516, 2).PasteSpecial(-4122)
$ws.Cells.Item(516, 2).Value = 'CoCrFeNi'
$ws.Cells.Item(516, 3).Value = 'FCC'
$ws.Cells.Item(516, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(516, 4).Copy()
$ws.Cells.Item(516, 5).PasteSpecial(-4122)
$ws.Cells.Item(516, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(516, 6).Value = 'tensile yield stress'
$ws.Cells.Item(516, 7).Value = 'EXP'
$ws.Cells.Item(516, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(516, 9).Value = 298
$ws.Cells.Item(516, 10).Value = 209000000
$ws.Cells.Item(516, 11).Value = 8000000
$ws.Cells.Item(516, 12).Value = 'Pa'
$ws.Cells.Item(516, 13).Value = 'T2'
$ws.Cells.Item(516, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 517: CoCrFeNiGe0.1 / tensile yield stress
$ws.Cells.Item(517, 2).Value = 'CoCrFeNiGe0.1'
$ws.Cells.Item(517, 3).Value = 'FCC'
$ws.Cells.Item(517, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(517, 4).Copy()
$ws.Cells.Item(517, 5).PasteSpecial(-4122)
$ws.Cells.Item(517, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(517, 6).Value = 'tensile yield stress'
$ws.Cells.Item(517, 7).Value = 'EXP'
$ws.Cells.Item(517, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(517, 9).Value = 298
$ws.Cells.Item(517, 10).Value = 213000000
$ws.Cells.Item(517, 11).Value = 9000000
$ws.Cells.Item(517, 12).Value = 'Pa'
$ws.Cells.Item(517, 13).Value = 'T2'
$ws.Cells.Item(517, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 518: CoCrFeNiGe0.2 / tensile yield stress
$ws.Cells.Item(518, 2).Value = 'CoCrFeNiGe0.2'
$ws.Cells.Item(518, 3).Value = 'FCC'
$ws.Cells.Item(518, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(518, 4).Copy()
$ws.Cells.Item(518, 5).PasteSpecial(-4122)
$ws.Cells.Item(518, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(518, 6).Value = 'tensile yield stress'
$ws.Cells.Item(518, 7).Value = 'EXP'
$ws.Cells.Item(518, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(518, 9).Value = 298
$ws.Cells.Item(518, 10).Value = 218000000
$ws.Cells.Item(518, 11).Value = 12000000
$ws.Cells.Item(518, 12).Value = 'Pa'
$ws.Cells.Item(518, 13).Value = 'T2'
$ws.Cells.Item(518, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 519: CoCrFeNiGe0.3 / tensile yield stress
$ws.Cells.Item(519, 2).Value = 'CoCrFeNiGe0.3'
$ws.Cells.Item(519, 3).Value = 'FCC'
$ws.Cells.Item(519, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(519, 4).Copy()
$ws.Cells.Item(519, 5).PasteSpecial(-4122)
$ws.Cells.Item(519, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(519, 6).Value = 'tensile yield stress'
$ws.Cells.Item(519, 7).Value = 'EXP'
$ws.Cells.Item(519, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(519, 9).Value = 298
$ws.Cells.Item(519, 10).Value = 223000000
$ws.Cells.Item(519, 11).Value = 3000000
$ws.Cells.Item(519, 12).Value = 'Pa'
$ws.Cells.Item(519, 13).Value = 'T2'
$ws.Cells.Item(519, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 520: CoCrFeNi / UTS
$ws.Cells.Item(488, 2).Copy()
$ws.Cells.Item(520, 2).PasteSpecial(-4122)
$ws.Cells.Item(520, 2).Value = 'CoCrFeNi'
$ws.Cells.Item(520, 3).Value = 'FCC'
$ws.Cells.Item(520, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(520, 4).Copy()
$ws.Cells.Item(520, 5).PasteSpecial(-4122)
$ws.Cells.Item(520, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(520, 6).Value = 'UTS'
$ws.Cells.Item(520, 7).Value = 'EXP'
$ws.Cells.Item(520, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(520, 9).Value = 298
$ws.Cells.Item(520, 10).Value = 590000000
$ws.Cells.Item(520, 11).Value = 16000000
$ws.Cells.Item(520, 12).Value = 'Pa'
$ws.Cells.Item(520, 13).Value = 'T2'
$ws.Cells.Item(520, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 521: CoCrFeNiGe0.1 / UTS
$ws.Cells.Item(521, 2).Value = 'CoCrFeNiGe0.1'
$ws.Cells.Item(521, 3).Value = 'FCC'
$ws.Cells.Item(521, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(521, 4).Copy()
$ws.Cells.Item(521, 5).PasteSpecial(-4122)
$ws.Cells.Item(521, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(521, 6).Value = 'UTS'
$ws.Cells.Item(521, 7).Value = 'EXP'
$ws.Cells.Item(521, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(521, 9).Value = 298
$ws.Cells.Item(521, 10).Value = 592000000
$ws.Cells.Item(521, 11).Value = 20000000
$ws.Cells.Item(521, 12).Value = 'Pa'
$ws.Cells.Item(521, 13).Value = 'T2'
$ws.Cells.Item(521, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 522: CoCrFeNiGe0.2 / UTS
$ws.Cells.Item(522, 2).Value = 'CoCrFeNiGe0.2'
$ws.Cells.Item(522, 3).Value = 'FCC'
$ws.Cells.Item(522, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(522, 4).Copy()
$ws.Cells.Item(522, 5).PasteSpecial(-4122)
$ws.Cells.Item(522, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(522, 6).Value = 'UTS'
$ws.Cells.Item(522, 7).Value = 'EXP'
$ws.Cells.Item(522, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(522, 9).Value = 298
$ws.Cells.Item(522, 10).Value = 596000000
$ws.Cells.Item(522, 11).Value = 13000000
$ws.Cells.Item(522, 12).Value = 'Pa'
$ws.Cells.Item(522, 13).Value = 'T2'
$ws.Cells.Item(522, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 523: CoCrFeNiGe0.3 / UTS
$ws.Cells.Item(523, 2).Value = 'CoCrFeNiGe0.3'
$ws.Cells.Item(523, 3).Value = 'FCC'
$ws.Cells.Item(523, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(523, 4).Copy()
$ws.Cells.Item(523, 5).PasteSpecial(-4122)
$ws.Cells.Item(523, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(523, 6).Value = 'UTS'
$ws.Cells.Item(523, 7).Value = 'EXP'
$ws.Cells.Item(523, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(523, 9).Value = 298
$ws.Cells.Item(523, 10).Value = 617000000
$ws.Cells.Item(523, 11).Value = 17000000
$ws.Cells.Item(523, 12).Value = 'Pa'
$ws.Cells.Item(523, 13).Value = 'T2'
$ws.Cells.Item(523, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 524: CoCrFeNi / tensile ductility
$ws.Cells.Item(488, 2).Copy()
$ws.Cells.Item(524, 2).PasteSpecial(-4122)
$ws.Cells.Item(524, 2).Value = 'CoCrFeNi'
$ws.Cells.Item(524, 3).Value = 'FCC'
$ws.Cells.Item(524, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(524, 4).Copy()
$ws.Cells.Item(524, 5).PasteSpecial(-4122)
$ws.Cells.Item(524, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(524, 6).Value = 'tensile ductility'
$ws.Cells.Item(524, 7).Value = 'EXP'
$ws.Cells.Item(524, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(524, 9).Value = 298
$ws.Cells.Item(524, 10).Value = 60
$ws.Cells.Item(524, 11).Value = 3.2
$ws.Cells.Item(524, 12).Value = '%'
$ws.Cells.Item(524, 13).Value = 'T2'
$ws.Cells.Item(524, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 525: CoCrFeNiGe0.1 / tensile ductility
$ws.Cells.Item(525, 2).Value = 'CoCrFeNiGe0.1'
$ws.Cells.Item(525, 3).Value = 'FCC'
$ws.Cells.Item(525, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(525, 4).Copy()
$ws.Cells.Item(525, 5).PasteSpecial(-4122)
$ws.Cells.Item(525, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(525, 6).Value = 'tensile ductility'
$ws.Cells.Item(525, 7).Value = 'EXP'
$ws.Cells.Item(525, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(525, 9).Value = 298
$ws.Cells.Item(525, 10).Value = 60.2
$ws.Cells.Item(525, 11).Value = 2.5
$ws.Cells.Item(525, 12).Value = '%'
$ws.Cells.Item(525, 13).Value = 'T2'
$ws.Cells.Item(525, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 526: CoCrFeNiGe0.2 / tensile ductility
$ws.Cells.Item(526, 2).Value = 'CoCrFeNiGe0.2'
$ws.Cells.Item(526, 3).Value = 'FCC'
$ws.Cells.Item(526, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(526, 4).Copy()
$ws.Cells.Item(526, 5).PasteSpecial(-4122)
$ws.Cells.Item(526, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(526, 6).Value = 'tensile ductility'
$ws.Cells.Item(526, 7).Value = 'EXP'
$ws.Cells.Item(526, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(526, 9).Value = 298
$ws.Cells.Item(526, 10).Value = 62.8
$ws.Cells.Item(526, 11).Value = 2.1
$ws.Cells.Item(526, 12).Value = '%'
$ws.Cells.Item(526, 13).Value = 'T2'
$ws.Cells.Item(526, 14).Value = '10.1016/j.intermet.2021.107167'

# Row 527: CoCrFeNiGe0.3 / tensile ductility
$ws.Cells.Item(527, 2).Value = 'CoCrFeNiGe0.3'
$ws.Cells.Item(527, 3).Value = 'FCC'
$ws.Cells.Item(527, 4).Value = 'VAM+H+WQ'
$ws.Cells.Item(527, 4).Copy()
$ws.Cells.Item(527, 5).PasteSpecial(-4122)
$ws.Cells.Item(527, 5).Value = 'homogenized at 1273K for 24h and water quenched'
$ws.Cells.Item(527, 6).Value = 'tensile ductility'
$ws.Cells.Item(527, 7).Value = 'EXP'
$ws.Cells.Item(527, 8).Value = 'strain rate 1e-3/s'
$ws.Cells.Item(527, 9).Value = 298
$ws.Cells.Item(527, 10).Value = 63.2
$ws.Cells.Item(527, 11).Value = 4.5
$ws.Cells.Item(527, 12).Value = '%'
$ws.Cells.Item(527, 13).Value = 'T2'
$ws.Cells.Item(527, 14).Value = '10.1016/j.intermet.2021.107167'

$excel.CutCopyMode = 0

# --- Step 3: scroll/selection state matching the post-edit view ---
$ws.Application.ActiveWindow.ScrollRow = 481
$ws.Range("K535").Select()